# Auto update stock data
# Update the "Date_1" values in column A from 2026/01/01 to 2026/01/02
# for every ticker block's first (most recent) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $targetRows) {
    $cell = $ws.Range("A" + $r)

    # Force a text format before assigning so the slash-delimited string
    # isn't auto-parsed into a date serial number by the COM layer, then
    # reset the style back to Normal so no stray number-format/style is
    # left attached to the cell (keeps it identical to the original,
    # un-styled inline-string cell other than the text itself).
    $cell.NumberFormat = "@"
    $cell.Value = "2026/01/02"
    $cell.Style = "Normal"
}
